# Update scripts with new TPM (transcripts per million) values.
#
# The ligand "Ligand average expression value" (G) / "Ligand total expression
# value" (H) columns were recomputed from the new TPM data for the three
# sending clusters (ECs rows 2-3, FAPs rows 4-5, MuSCs rows 6-7). All of the
# specificity / edge-weight columns that are mathematically derived from G
# and H (I, J, Q, R, S, T) were recomputed downstream and are updated here to
# match, cell for cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (ECs -> FAPs) ----
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("Q2").Value = 179.983350602076
$ws.Range("R2").Value = 1619.850155418684
$ws.Range("S2").Value = 0.4449836184688483
$ws.Range("T2").Value = 0.4449836184688483

# ---- Row 3 (ECs -> MuSCs) ----
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("Q3").Value = 1.604456161893667
$ws.Range("R3").Value = 14.440105457043
$ws.Range("S3").Value = 0.003966793073946968
$ws.Range("T3").Value = 0.003966793073946969

# ---- Row 4 (FAPs -> FAPs) ----
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("S4").Value = 0.3385626345409294
$ws.Range("T4").Value = 0.3385626345409294

# ---- Row 5 (FAPs -> MuSCs) ----
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("S5").Value = 0.003018106415726889
$ws.Range("T5").Value = 0.00301810641572689

# ---- Row 6 (MuSCs -> FAPs) ----
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("Q6").Value = 83.97565532983201
$ws.Range("R6").Value = 755.780897968488
$ws.Range("S6").Value = 0.2076180427076151
$ws.Range("T6").Value = 0.2076180427076151

# ---- Row 7 (MuSCs -> MuSCs) ----
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("Q7").Value = 0.7485984519806667
$ws.Range("S7").Value = 0.001850804792933403
$ws.Range("T7").Value = 0.001850804792933404
